# ManagePlaylist task: Fetch by Artist
# Insert a new bulleted list item (with its "Multiplication Sign" icon)
# right after the "...which identifies text picked, (hidden field)" bullet
# in the Fetch(15)/ArtistArg(1) events cell, directly before the
# "ArtistArg id to SearchArgID; (hidden field)" bullet.

$d = $word.ActiveDocument

$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "which identifies text picked, (hidden field)",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if (-not $found) {
    throw "Anchor paragraph for the new 'Validate artists partial name is present' bullet was not found."
}

# Collapse to the paragraph that holds the match, then move to the very end
# of that paragraph (i.e. just after its paragraph mark / right before the
# next bullet's first run) so the new content lands as its own <w:p>.
$anchorParagraph = $searchRange.Paragraphs(1)
$insertPoint = $anchorParagraph.Range.End
$insertionRange = $d.Range($insertPoint, $insertPoint)

$newParagraphXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" mc:Ignorable="w14 wp14"><w:body><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/><w:rPr><w:lang w:val="en-CA"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:val="en-CA"/></w:rPr><mc:AlternateContent><mc:Choice Requires="wps"><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251671552" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="51B956D6" wp14:editId="333A9B6A"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>-46714</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>145415</wp:posOffset></wp:positionV><wp:extent cx="914400" cy="914400"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapNone/><wp:docPr id="7" name="Multiplication Sign 7"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvSpPr/><wps:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="914400" cy="914400"/></a:xfrm><a:prstGeom prst="mathMultiply"><a:avLst/></a:prstGeom></wps:spPr><wps:style><a:lnRef idx="2"><a:schemeClr val="accent1"><a:shade val="50000"/></a:schemeClr></a:lnRef><a:fillRef idx="1"><a:schemeClr val="accent1"/></a:fillRef><a:effectRef idx="0"><a:schemeClr val="accent1"/></a:effectRef><a:fontRef idx="minor"><a:schemeClr val="lt1"/></a:fontRef></wps:style><wps:bodyPr rot="0" spcFirstLastPara="0" vertOverflow="overflow" horzOverflow="overflow" vert="horz" wrap="square" lIns="91440" tIns="45720" rIns="91440" bIns="45720" numCol="1" spcCol="0" rtlCol="0" fromWordArt="0" anchor="ctr" anchorCtr="0" forceAA="0" compatLnSpc="1"><a:prstTxWarp prst="textNoShape"><a:avLst/></a:prstTxWarp><a:noAutofit/></wps:bodyPr></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback><w:pict><v:shape w14:anchorId="2B3466B1" id="Multiplication Sign 7" o:spid="_x0000_s1026" style="position:absolute;margin-left:-3.7pt;margin-top:11.45pt;width:1in;height:1in;z-index:251671552;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;v-text-anchor:middle" coordsize="914400,914400" o:gfxdata="UEsDBBQABgAIAAAAIQC2gziS/gAAAOEBAAATAAAAW0NvbnRlbnRfVHlwZXNdLnhtbJSRQU7DMBBF&#10;90jcwfIWJU67QAgl6YK0S0CoHGBkTxKLZGx5TGhvj5O2G0SRWNoz/78nu9wcxkFMGNg6quQqL6RA&#10;0s5Y6ir5vt9lD1JwBDIwOMJKHpHlpr69KfdHjyxSmriSfYz+USnWPY7AufNIadK6MEJMx9ApD/oD&#10;OlTrorhX2lFEilmcO2RdNtjC5xDF9pCuTyYBB5bi6bQ4syoJ3g9WQ0ymaiLzg5KdCXlKLjvcW893&#10;SUOqXwnz5DrgnHtJTxOsQfEKIT7DmDSUCaxw7Rqn8787ZsmRM9e2VmPeBN4uqYvTtW7jvijg9N/y&#10;JsXecLq0q+WD6m8AAAD//wMAUEsDBBQABgAIAAAAIQA4/SH/1gAAAJQBAAALAAAAX3JlbHMvLnJl&#10;bHOkkMFqwzAMhu+DvYPRfXGawxijTi+j0GvpHsDYimMaW0Yy2fr2M4PBMnrbUb/Q94l/f/hMi1qR&#10;JVI2sOt6UJgd+ZiDgffL8ekFlFSbvV0oo4EbChzGx4f9GRdb25HMsYhqlCwG5lrLq9biZkxWOiqY&#10;22YiTra2kYMu1l1tQD30/bPm3wwYN0x18gb45AdQl1tp5j/sFB2T0FQ7R0nTNEV3j6o9feQzro1i&#10;OWA14Fm+Q8a1a8+Bvu/d/dMb2JY5uiPbhG/ktn4cqGU/er3pcvwCAAD//wMAUEsDBBQABgAIAAAA&#10;IQC7yBCzeAIAAEoFAAAOAAAAZHJzL2Uyb0RvYy54bWysVE1vGjEQvVfqf7B8b3ZBpLQoS4SIUlWi&#10;CQqpcna8NmvJ9ri2YaG/vmPvskFJ1ENVDma8M/Pm642vrg9Gk73wQYGt6OiipERYDrWy24r+fLz9&#10;9IWSEJmtmQYrKnoUgV7PP364at1MjKEBXQtPEMSGWesq2sToZkUReCMMCxfghEWlBG9YxKvfFrVn&#10;LaIbXYzL8nPRgq+dBy5CwK83nZLOM76Ugsd7KYOIRFcUc4v59Pl8Tmcxv2KzrWeuUbxPg/1DFoYp&#10;i0EHqBsWGdl59QbKKO4hgIwXHEwBUioucg1Yzah8Vc2mYU7kWrA5wQ1tCv8Plt/t156ouqJTSiwz&#10;OKIfOx2V04qziEMlG7W1ZJoa1bowQ/uNW/v+FlBMVR+kN+kf6yGH3Nzj0FxxiITjx6+jyaTEEXBU&#10;9TKiFC/Ozof4TYAhSagoDrzpUznm1rL9KsTO5WSK/impLo0sxaMWKRNtH4TEujDwOHtnRoml9mTP&#10;kAuMc2HjqFM1rBbd58sSf6lWzGvwyLcMmJCl0nrA7gESW99idzC9fXIVmZCDc/m3xDrnwSNHBhsH&#10;Z6Ms+PcANFbVR+7sT03qWpO69Az1EafuoVuH4Pitwp6vWIhr5pH/OCbc6XiPh9TQVhR6iZIG/O/3&#10;vid7pCVqKWlxnyoafu2YF5To7xYJm0eOC5gvk8vpGGP4c83zucbuzBJwTCN8PRzPYrKP+iRKD+YJ&#10;V3+RoqKKWY6xK8qjP12WsdtzfDy4WCyyGS6dY3FlN44n8NTVxKXHwxPzrideRMbewWn32OwV7zrb&#10;5GlhsYsgVSblS1/7fuPCZuL0j0t6Ec7v2erlCZz/AQAA//8DAFBLAwQUAAYACAAAACEAIZv8cOIA&#10;AAAJAQAADwAAAGRycy9kb3ducmV2LnhtbEyPwU7DMBBE70j8g7VIXFDrkCJDQ5yqqoQQB1TRVurV&#10;jZckJF4H22lDvx73BLdZzWjmbb4YTceO6HxjScL9NAGGVFrdUCVht32ZPAHzQZFWnSWU8IMeFsX1&#10;Va4ybU/0gcdNqFgsIZ8pCXUIfca5L2s0yk9tjxS9T+uMCvF0FddOnWK56XiaJIIb1VBcqFWPqxrL&#10;djMYCec7VQ7uNcxas37/3r+126/l6izl7c24fAYWcAx/YbjgR3QoItPBDqQ96yRMHh9iUkKazoFd&#10;/JkQwA5RCDEHXuT8/wfFLwAAAP//AwBQSwECLQAUAAYACAAAACEAtoM4kv4AAADhAQAAEwAAAAAA&#10;AAAAAAAAAAAAAAAAW0NvbnRlbnRfVHlwZXNdLnhtbFBLAQItABQABgAIAAAAIQA4/SH/1gAAAJQB&#10;AAALAAAAAAAAAAAAAAAAAC8BAABfcmVscy8ucmVsc1BLAQItABQABgAIAAAAIQC7yBCzeAIAAEoF&#10;AAAOAAAAAAAAAAAAAAAAAC4CAABkcnMvZTJvRG9jLnhtbFBLAQItABQABgAIAAAAIQAhm/xw4gAA&#10;AAkBAAAPAAAAAAAAAAAAAAAAANIEAABkcnMvZG93bnJldi54bWxQSwUGAAAAAAQABADzAAAA4QUA&#10;AAAA&#10;" path="m143578,295654l295654,143578,457200,305125,618746,143578,770822,295654,609275,457200,770822,618746,618746,770822,457200,609275,295654,770822,143578,618746,305125,457200,143578,295654xe" fillcolor="#5b9bd5 [3204]" strokecolor="#1f4d78 [1604]" strokeweight="1pt"><v:stroke joinstyle="miter"/><v:path arrowok="t" o:connecttype="custom" o:connectlocs="143578,295654;295654,143578;457200,305125;618746,143578;770822,295654;609275,457200;770822,618746;618746,770822;457200,609275;295654,770822;143578,618746;305125,457200;143578,295654" o:connectangles="0,0,0,0,0,0,0,0,0,0,0,0,0"/></v:shape></w:pict></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:lang w:val="en-CA"/></w:rPr><w:t>Validate artists partial name is present</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionRange.InsertXML($newParagraphXml)
